# Re-added page object fixes, mega menu tests
#
# 1. Add a new worksheet "MegaMenuInfo" right after "HomeLandingTopicCards".
# 2. Populate it with the mega-menu test fixture data.
# 3. Copy the header-row formatting (bold + shaded fill) from the existing
#    sheet so no new styles are introduced.
# 4. Reasonable column widths for the new sheet.
# 5. Adjust selection/active-sheet bookkeeping to match: the first sheet's
#    selection moves to A1:C1 (and it is no longer the active tab), and the
#    new sheet becomes the active tab with the cursor left under the data.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Create the new sheet right after the existing one ---------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "MegaMenuInfo"

# --- Fill in the data (order chosen to match the shared-string table) ------
$ws2.Range("E1").Value = "Language"
$ws2.Range("E2").Value = "english"
$ws2.Range("E3").Value = "spanish"

$ws2.Range("B1").Value = "NavGroup"
$ws2.Range("C1").Value = "SubNavGroup"
$ws2.Range("D1").Value = "ListItem"

$ws2.Range("B3").Value = "Tipos de cáncer"
$ws2.Range("C3").Value = "Tipos comunes de cáncer"
$ws2.Range("D3").Value = "Linfoma"

$ws2.Range("B2").Value = "About Cancer"
$ws2.Range("C2").Value = "Understanding Cancer"
$ws2.Range("D2").Value = "Cancer Statistics"

$ws2.Range("A3").Value = "/espanol/investigacion"
$ws2.Range("A1").Value = "Path"
$ws2.Range("A2").Value = "/"

# --- Header formatting: reuse the existing header style (bold + fill) ------
# Source A1 uses the plain bold+fill header style (others in row 1 mix in a
# text-numberformat variant), so copy that single cell's format across the
# whole new header row.
# -4122 == xlPasteFormats
$ws1.Range("A1").Copy() | Out-Null
$ws2.Range("A1:E1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Column widths on the new sheet -----------------------------------------
$ws2.Columns.Item(1).ColumnWidth = 20.6
$ws2.Columns.Item(2).ColumnWidth = 23.3
$ws2.Columns.Item(3).ColumnWidth = 22.6
$ws2.Columns.Item(4).ColumnWidth = 23.3
$ws2.Columns.Item(5).ColumnWidth = 23.3

# --- Selection bookkeeping ---------------------------------------------------
# First sheet: selection becomes A1:C1, and it stops being the active tab.
$ws1.Range("A1:C1").Select() | Out-Null

# New sheet becomes active, cursor parked just below the pasted data.
$ws2.Activate() | Out-Null
$ws2.Range("A4").Select() | Out-Null

Write-Host "Added MegaMenuInfo worksheet with mega-menu test fixture data"
